$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 (Ano 2025) with refreshed figures
$ws.Range("B9").Value = 4077135.58
$ws.Range("C9").Value = 640931.24
$ws.Range("D9").Value = 4718066.82
$ws.Range("E9").Value = 13.58461557354544
$ws.Range("F9").Value = 86.41538442645455
$ws.Range("G9").Value = -38.05719148449573
$ws.Range("H9").Value = -26.37256700351249
$ws.Range("I9").Value = 41025
$ws.Range("J9").Value = 1759
$ws.Range("K9").Value = 42784
$ws.Range("L9").Value = 29638
$ws.Range("M9").Value = 159.1897840610028
$ws.Range("N9").Value = 8.682107912893144
